$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and volume-change (E) columns with the latest scrape.
$ws.Range("D2").Value = '20.209.58'
$ws.Range("E2").Value = '  -0.77%  '
$ws.Range("D3").Value = '1.431.43'
$ws.Range("E3").Value = '  -0.54%  '
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").Value = "'0.9948"
$ws.Range("E5").Value = '  -0.70%  '
$ws.Range("D6").Value = "'277.45"
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").Value = "'0.3709"
$ws.Range("E7").Value = '  -0.58%  '
$ws.Range("D8").Value = "'0.3160"
$ws.Range("E8").Value = '  +2.25%  '
$ws.Range("D9").Value = "'40.36"
$ws.Range("E9").Value = '  -0.63%  '
$ws.Range("D10").Value = "'1.062"
$ws.Range("E10").Value = '  +4.79%  '
$ws.Range("D11").Value = "'0.06596"
$ws.Range("E11").Value = '  -0.07%  '
$ws.Range("D12").Value = "'0.9976"
$ws.Range("E12").Value = '  -0.53%  '
$ws.Range("D13").Value = "'5.571"
$ws.Range("E13").Value = '  +3.66%  '
$ws.Range("D14").Value = "'18.30"
$ws.Range("E14").Value = '  +5.61%  '
$ws.Range("D15").Value = "'6.241"
$ws.Range("E15").Value = '  +1.49%  '
$ws.Range("D16").Value = "'0.00001032"
$ws.Range("E16").Value = '  +2.06%  '
$ws.Range("D17").Value = '1.432.13'
$ws.Range("E17").Value = '  -0.63%  '
$ws.Range("D18").Value = "'0.05759"
$ws.Range("E18").Value = '  -1.79%  '
$ws.Range("D19").Value = "'0.9943"
$ws.Range("E19").Value = '  -0.72%  '
$ws.Range("D20").Value = "'71.86"
$ws.Range("E20").Value = '  -6.21%  '
$ws.Range("D21").Value = "'5.635"
$ws.Range("E21").Value = '  -1.76%  '
$ws.Range("E22").Value = '  +3.47%  '
$ws.Range("D23").Value = "'11.15"
$ws.Range("E23").Value = '  +1.46%  '
$ws.Range("D24").Value = "'2.230"
$ws.Range("E24").Value = '  -4.22%  '
$ws.Range("D25").Value = '20.230.74'
$ws.Range("E25").Value = '  -0.63%  '
$ws.Range("D26").Value = "'2.319"
$ws.Range("E26").Value = '  +3.06%  '
$ws.Range("D27").Value = "'135.36"
$ws.Range("E27").Value = '  -5.22%  '
$ws.Range("E28").Value = '  +2.40%  '
$ws.Range("D29").Value = '1.593.70'
$ws.Range("E29").Value = '  -0.57%  '
$ws.Range("D30").Value = "'111.91"
$ws.Range("E30").Value = '  +1.56%  '
$ws.Range("D31").Value = "'3.961"
$ws.Range("E31").Value = '  -0.33%  '
$ws.Range("D32").Value = "'5.338"
$ws.Range("E32").Value = '  -2.87%  '
$ws.Range("D33").Value = "'0.8460"
$ws.Range("E33").Value = '  -8.49%  '
$ws.Range("D34").Value = "'0.07814"
$ws.Range("E34").Value = '  +1.22%  '
$ws.Range("D35").Value = "'1.497"
$ws.Range("E35").Value = '  +11.58%  '
$ws.Range("D36").Value = "'0.05932"
$ws.Range("E36").Value = '  +3.50%  '
$ws.Range("D37").Value = "'4.940"
$ws.Range("E37").Value = '  +3.96%  '
$ws.Range("D38").Value = "'10.85"
$ws.Range("E38").Value = '  -1.34%  '
$ws.Range("D39").Value = "'0.9949"
$ws.Range("E39").Value = '  -0.64%  '
$ws.Range("D40").Value = "'7.824"
$ws.Range("E40").Value = '  -6.26%  '
$ws.Range("D41").Value = "'0.02069"
$ws.Range("E41").Value = '  +1.92%  '
$ws.Range("D42").Value = "'1.115"
$ws.Range("E42").Value = '  -1.90%  '
$ws.Range("D43").Value = "'0.1888"
$ws.Range("E43").Value = '  -2.06%  '
$ws.Range("D44").Value = "'0.5388"
$ws.Range("E44").Value = '  +0.77%  '
$ws.Range("D45").Value = "'12.50"
$ws.Range("E45").Value = '  +3.49%  '
$ws.Range("D46").Value = "'3.564"
$ws.Range("E46").Value = '  -0.77%  '
$ws.Range("E47").Value = '  +6.37%  '
$ws.Range("D48").Value = "'0.5301"
$ws.Range("E48").Value = '  +2.37%  '
$ws.Range("D49").Value = "'1.809"
$ws.Range("E49").Value = '  +1.16%  '
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("D51").Value = "'0.06289"
$ws.Range("E51").Value = '  +0.26%  '
